$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Add two new vulnerability rows at the bottom of the list (rows 38-39)
$ws.Range("A38").Value = "Network"
$ws.Range("B38").Value = "p,a,a+"
$ws.Range("A39").Value = "Node"
$ws.Range("B39").Value = "p,a,a+"

$ws.Range("C38").Value = "PTV-NET-MITM-ICMP6REDIR"
$ws.Range("C39").Value = "PTV-NET-MITM-ICMP6REDIRDEV"

$ws.Range("D38").Value = "Network does not block ICMPv6 Redirect messages"
$ws.Range("D39").Value = "Device communication can be redirected using ICMPv6 Redirect"

# Move selection, mirroring the author's final cursor position
$ws.Range("K36").Select()
